$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for columns F ("height") and G ("weight") ---
$ws.Range("F1").Value = "height"
$ws.Range("G1").Value = "weight"

# Copy the header formatting (bold, centered, bordered) from the existing
# "fantasy points" header (E1) onto the two new header cells.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# --- Move the old "fantasy points" data (column E, rows 2-16) into the new
#     "weight" column (G), preserving the original numeric values exactly. ---
$ws.Range("E2:E16").Cut($ws.Range("G2:G16"))

# --- Populate the new constant columns: E = 6.5 (rewritten "fantasy points"
#     placeholder) and F = 255 ("height") for every data row. ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5
    $ws.Cells.Item($r, 6).Value = 255
}
